# Wed May  8 12:22:21 UTC 2024 - cryptos list refresh (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "584.84", "62.329.26") need an explicit text format so the literal
# string from the crypto-price feed is preserved, matching the source diff.
$textCells = 'D5','D6','D8','D11','D14','D18','D20','D21','D23','D24','D26','D30','D32','D33','D39','D41','D44','D45','D51'
foreach ($ref in $textCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range('D2').Value = '62.329.26'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '3.000.04'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '584.84'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').Value = '145.79'
$ws.Range('E6').Value = '  -6.08%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').Value = '3.000.07'
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('E10').Value = '  -6.24%  '
$ws.Range('D11').Value = '5.76'
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('E12').Value = '  -2.41%  '
$ws.Range('E13').Value = '  -4.40%  '
$ws.Range('D14').Value = '34.42'
$ws.Range('E14').Value = '  -5.84%  '
$ws.Range('E15').Value = '  +2.54%  '
$ws.Range('D16').Value = '3.498.61'
$ws.Range('E16').Value = '  -2.16%  '
$ws.Range('D17').Value = '62.363.65'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '6.97'
$ws.Range('E18').Value = '  -3.02%  '
$ws.Range('D19').Value = '3.002.03'
$ws.Range('E19').Value = '  -2.26%  '
$ws.Range('D20').Value = '454.01'
$ws.Range('E20').Value = '  -5.64%  '
$ws.Range('D21').Value = '13.79'
$ws.Range('E21').Value = '  -4.20%  '
$ws.Range('E22').Value = '  -4.32%  '
$ws.Range('D23').Value = '7.38'
$ws.Range('E23').Value = '  -1.97%  '
$ws.Range('D24').Value = '80.51'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('E25').Value = '  -6.75%  '
$ws.Range('D26').Value = '12.20'
$ws.Range('E26').Value = '  -4.77%  '
$ws.Range('D27').Value = '10.09'
$ws.Range('E27').Value = '  -4.67%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '7.17'
$ws.Range('E30').Value = '  -5.51%  '
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').Value = '2.09'
$ws.Range('E32').Value = '  -4.37%  '
$ws.Range('D33').Value = '26.83'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('E34').Value = '  -3.57%  '
$ws.Range('E35').Value = '  -3.90%  '
$ws.Range('D36').Value = '0.0₃0786'
$ws.Range('E36').Value = '  -5.91%  '
$ws.Range('E37').Value = '  -4.63%  '
$ws.Range('E38').Value = '  -5.77%  '
$ws.Range('D39').Value = '50.02'
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  -3.30%  '
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  -12.60%  '
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('D43').Value = '382.87'
$ws.Range('E43').Value = '  -12.54%  '
$ws.Range('D44').Value = '0.271'
$ws.Range('E44').Value = '  -6.64%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0352'
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.755.42'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('D47').Value = '38.43'
$ws.Range('E47').Value = '  -4.04%  '
$ws.Range('D48').Value = '127.95'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('D51').Value = '23.79'
$ws.Range('E51').Value = '  -6.48%  '
